$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value: beach frequency count changes from 18 to 14
$ws.Range("B3").Value = 14

# Add new row: "back" search word with frequency count 2
$ws.Range("A6").Value = "back"
$ws.Range("B6").Value = 2
